$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update dimension / ensure sheet has enough rows; values below will naturally extend used range

$ws.Cells.Item(11, 1).Value = 'ACORN Participating Countries'
$ws.Cells.Item(11, 2).Value = 'Các nước tham gia vào nghiên cứu ACORN'
$ws.Cells.Item(11, 3).ClearContents()
$ws.Cells.Item(12, 1).Value = 'All ''orgname'' are provided.'
$ws.Cells.Item(12, 2).Value = 'Tất cả "orgname" được cung cấp'
$ws.Cells.Item(12, 3).ClearContents()
$ws.Cells.Item(13, 1).Value = 'All ''patid'' are provided.'
$ws.Cells.Item(13, 2).Value = 'Tất cả "patid" được cung cấp'
$ws.Cells.Item(13, 3).ClearContents()
$ws.Cells.Item(14, 1).Value = 'All ''specdate'' are provided.'
$ws.Cells.Item(14, 2).Value = 'Tất cả "specdate" được cung cấp'
$ws.Cells.Item(14, 3).ClearContents()
$ws.Cells.Item(15, 1).Value = 'All ''specdate'' are today or before today.'
$ws.Cells.Item(15, 2).Value = 'Tât cả "specdate" là ngày hôm này hoặc trước ngày hôm nay'
$ws.Cells.Item(15, 3).ClearContents()
$ws.Cells.Item(16, 1).Value = 'All ''specgroup'' are provided.'
$ws.Cells.Item(16, 2).Value = 'Tất cả "specgroup" được cung cấp'
$ws.Cells.Item(16, 3).ClearContents()
$ws.Cells.Item(17, 1).Value = 'All ''specid'' are provided.'
$ws.Cells.Item(17, 2).Value = 'Tất cả "specid" được cung cấp'
$ws.Cells.Item(17, 3).ClearContents()
$ws.Cells.Item(18, 1).Value = 'All dates of enrolment for HAI patients have a matching date in the HAI survey dataset'
$ws.Cells.Item(18, 2).Value = 'Tất cả các ngày thu tuyển của bệnh nhân HAI phải tương ứng với ngày thực hiện giám sát điểm đã được định sẵn'
$ws.Cells.Item(18, 3).ClearContents()
$ws.Cells.Item(19, 1).Value = 'All Other Organisms'
$ws.Cells.Item(19, 2).Value = 'Tất cả các vi sinh vật khác'
$ws.Cells.Item(19, 3).ClearContents()
$ws.Cells.Item(20, 1).Value = 'All valid records have an ACORN ID.'
$ws.Cells.Item(20, 2).Value = 'TBT'
$ws.Cells.Item(20, 3).Value = 'new'
$ws.Cells.Item(21, 1).Value = 'AMR'
$ws.Cells.Item(21, 2).Value = 'TBT'
$ws.Cells.Item(21, 3).Value = 'new'
$ws.Cells.Item(22, 1).Value = 'and generate enrolment log.'
$ws.Cells.Item(22, 2).Value = 'và tạo nhật ký thu tuyển.'
$ws.Cells.Item(22, 3).ClearContents()
$ws.Cells.Item(23, 1).Value = 'Attempting to connect.'
$ws.Cells.Item(23, 2).Value = 'Đang cố gắng kết nối.'
$ws.Cells.Item(23, 3).ClearContents()
$ws.Cells.Item(24, 1).Value = 'Blood culture collected within 24 hours of admission (CAI) / symptom onset (HAI)'
$ws.Cells.Item(24, 2).Value = 'Mẫu cấy máu thu thập trong vòng 24h sau khi nhập viện (CAI)/ khởi phát triệu chứng (HAI)'
$ws.Cells.Item(24, 3).ClearContents()
$ws.Cells.Item(25, 1).Value = 'Blood Culture Contaminants'
$ws.Cells.Item(25, 2).Value = 'Tạp nhiễm cấy máu'
$ws.Cells.Item(25, 3).ClearContents()
$ws.Cells.Item(26, 1).Value = 'Bloodstream Infection (BSI)'
$ws.Cells.Item(26, 2).Value = 'Nhiễm trùng huyết (BSI)'
$ws.Cells.Item(26, 3).ClearContents()
$ws.Cells.Item(27, 1).Value = 'Calculated age is consistent with ''Age Category'''
$ws.Cells.Item(27, 2).Value = 'Tuổi được tính toán nhất quán với "Phân loại tuổi"'
$ws.Cells.Item(27, 3).ClearContents()
$ws.Cells.Item(28, 1).Value = 'Calculated age isn''t always consistent with ''Age Category'''
$ws.Cells.Item(28, 2).Value = 'Tuổi được tính toán không nhất quán với "Phân loại tuổi"'
$ws.Cells.Item(28, 3).ClearContents()
$ws.Cells.Item(29, 1).Value = 'Cancel'
$ws.Cells.Item(29, 2).Value = 'Hủy bỏ'
$ws.Cells.Item(29, 3).ClearContents()
$ws.Cells.Item(30, 1).Value = 'Care should be taken when interpreting rates and AMR profiles where there are small numbers of cases or bacterial isolates: point estimates may be unreliable.'
$ws.Cells.Item(30, 2).Value = 'Cần thận trọng khi phiên giải về tỷ lệ và thông tin AMR khi có một số lượng nhỏ ca bệnh hoặc vi khuẩn được phân lập: ước tính điểm có thể không đáng tin cậy.'
$ws.Cells.Item(30, 3).ClearContents()
$ws.Cells.Item(31, 1).Value = 'Clinical and day-28 outcomes are consistent.'
$ws.Cells.Item(31, 2).Value = 'Kêt quả lâm sàng và ngày 28 là nhất quán'
$ws.Cells.Item(31, 3).ClearContents()
$ws.Cells.Item(32, 1).Value = 'Clinical and day-28 outcomes aren''t consistent for some dead patients.'
$ws.Cells.Item(32, 2).Value = 'Kết quả lâm sàng và ngày 28 là không nhất quán đối với một số bệnh nhân tử vong'
$ws.Cells.Item(32, 3).ClearContents()
$ws.Cells.Item(33, 1).Value = 'Clinical Outcome'
$ws.Cells.Item(33, 2).Value = 'Kết cục lâm sàng'
$ws.Cells.Item(33, 3).ClearContents()
$ws.Cells.Item(34, 1).Value = 'Clinical Outcome Status:'
$ws.Cells.Item(34, 2).Value = 'Tình trạng kết cục lâm sàng'
$ws.Cells.Item(34, 3).ClearContents()
$ws.Cells.Item(35, 1).Value = 'Co-resistances'
$ws.Cells.Item(35, 2).Value = 'TBT'
$ws.Cells.Item(35, 3).Value = 'new'
$ws.Cells.Item(36, 1).Value = 'Combine Susceptible + Intermediate'
$ws.Cells.Item(36, 2).Value = 'Kết hợp Nhạy cảm + Trung bình'
$ws.Cells.Item(36, 3).ClearContents()
$ws.Cells.Item(37, 1).Value = 'Consider saving .acorn file on the cloud for additional security.'
$ws.Cells.Item(37, 2).Value = 'Cân nhắc lưu tệp .acorn trên đám mây vì mục đích bảo mật.'
$ws.Cells.Item(37, 3).ClearContents()
$ws.Cells.Item(38, 1).Value = 'Contains names of organisms before and after mapping.'
$ws.Cells.Item(38, 2).Value = 'TBT'
$ws.Cells.Item(38, 3).Value = 'new'
$ws.Cells.Item(39, 1).Value = 'Couldn''t connect to server. Please check internet access.'
$ws.Cells.Item(39, 2).Value = 'Không thể kết nối với máy chủ. Vui lòng kiểm tra kết nối internet'
$ws.Cells.Item(39, 3).ClearContents()
$ws.Cells.Item(40, 1).Value = 'Critical errors with clinical data.'
$ws.Cells.Item(40, 2).Value = 'Lỗi nghiêm trọng với dữ liệu lâm sàng.'
$ws.Cells.Item(40, 3).ClearContents()
$ws.Cells.Item(41, 1).Value = 'Culture results per specimen type'
$ws.Cells.Item(41, 2).Value = 'Kết quả nuôi cấy cho mỗi loại mẫu'
$ws.Cells.Item(41, 3).ClearContents()
$ws.Cells.Item(42, 1).Value = 'Data Management'
$ws.Cells.Item(42, 2).Value = 'Quản lý dữ liệu'
$ws.Cells.Item(42, 3).ClearContents()
$ws.Cells.Item(43, 1).Value = 'Date of Enrolment'
$ws.Cells.Item(43, 2).Value = 'Ngày thu tuyển'
$ws.Cells.Item(43, 3).ClearContents()
$ws.Cells.Item(44, 1).Value = 'Day 28'
$ws.Cells.Item(44, 2).Value = 'Ngày 28'
$ws.Cells.Item(44, 3).ClearContents()
$ws.Cells.Item(45, 1).Value = 'Day 28 Status:'
$ws.Cells.Item(45, 2).Value = 'Tình trạng ngày 28'
$ws.Cells.Item(45, 3).ClearContents()
$ws.Cells.Item(46, 1).Value = 'Diagnosis at Enrolment'
$ws.Cells.Item(46, 2).Value = 'Chẩn đoán tại thời điểm thu tuyển'
$ws.Cells.Item(46, 3).ClearContents()
$ws.Cells.Item(47, 1).Value = 'Dismiss'
$ws.Cells.Item(47, 2).Value = 'TBT'
$ws.Cells.Item(47, 3).Value = 'new'
$ws.Cells.Item(48, 1).Value = 'Distribution of Enrolments'
$ws.Cells.Item(48, 2).Value = 'Phân bố thu tuyển'
$ws.Cells.Item(48, 3).ClearContents()
$ws.Cells.Item(49, 1).Value = 'Download Enrolment Log (.xlsx)'
$ws.Cells.Item(49, 2).Value = 'Tải xuống Sổ thu tuyển (.xlsx)'
$ws.Cells.Item(49, 3).ClearContents()
$ws.Cells.Item(50, 1).Value = 'Download Lab Log (.xlsx)'
$ws.Cells.Item(50, 2).Value = 'TBT'
$ws.Cells.Item(50, 3).Value = 'new'
$ws.Cells.Item(51, 1).Value = 'Empiric Antibiotics Prescribed'
$ws.Cells.Item(51, 2).Value = 'Điều trị kháng sinh theo kinh nghiệm'
$ws.Cells.Item(51, 3).ClearContents()
$ws.Cells.Item(52, 1).Value = 'Enrolments'
$ws.Cells.Item(52, 2).Value = 'Thu tuyển'
$ws.Cells.Item(52, 3).ClearContents()
$ws.Cells.Item(53, 1).Value = 'Enrolments by (type of) Ward'
$ws.Cells.Item(53, 2).Value = 'Thu tuyển theo (loại) Khoa'
$ws.Cells.Item(53, 3).ClearContents()
$ws.Cells.Item(54, 1).Value = 'Enrolments with Blood Culture'
$ws.Cells.Item(54, 2).Value = 'Thu tuyển với mẫu cấy máu'
$ws.Cells.Item(54, 3).ClearContents()
$ws.Cells.Item(55, 1).Value = 'Error in combining clinical and lab data.'
$ws.Cells.Item(55, 2).Value = 'Lỗi khi kết hợp dữ liệu lâm sàng và dữ liệu phòng xét nghiệm'
$ws.Cells.Item(55, 3).ClearContents()
$ws.Cells.Item(56, 1).Value = 'Every D28 record (F04) matches exactly one patient enrolment (F01).'
$ws.Cells.Item(56, 2).Value = 'Mỗi Phiếu dữ liệu ngày 28 (F04) khớp chính xác với một bệnh nhân thu tuyển (F01).'
$ws.Cells.Item(56, 3).ClearContents()
$ws.Cells.Item(57, 1).Value = 'Every hospital outcome record (F03) has a matching infection episode (F02).'
$ws.Cells.Item(57, 2).Value = 'Mỗi phiếu ra viện (F03) đều có phiếu nhiễm trùng (F02) tương ứng'
$ws.Cells.Item(57, 3).ClearContents()
$ws.Cells.Item(58, 1).Value = 'Every hospital outcome record (F03) has a matching patient enrolment (F01).'
$ws.Cells.Item(58, 2).Value = 'Mỗi phiếu xuất viện (F03) có một Phiếu thu tuyển (F01) tương ứng'
$ws.Cells.Item(58, 3).ClearContents()
$ws.Cells.Item(59, 1).Value = 'Every infection episode record (F02) has a matching patient enrolment (F01).'
$ws.Cells.Item(59, 2).Value = 'Mỗi Phiếu đợt nhiễm trùng (F02) có một phiếu thu tuyển tương ứng (F01)'
$ws.Cells.Item(59, 3).ClearContents()
$ws.Cells.Item(60, 1).Value = 'File name:'
$ws.Cells.Item(60, 2).Value = 'Tên tệp:'
$ws.Cells.Item(60, 3).ClearContents()
$ws.Cells.Item(61, 1).Value = 'First sheet is the log of all enrolments retrived from REDCap (as per adjacent table). The second sheet is a listing of all flagged elements.'
$ws.Cells.Item(61, 2).Value = 'Trang đầu tiên là nhật ký thu tuyển được truy xuất từ REDCap (theo bảng liền kề). Trang thứ hai là danh sách tất cả các yếu tố được gắn cờ cảnh báo.'
$ws.Cells.Item(61, 3).ClearContents()
$ws.Cells.Item(62, 1).Value = 'Follow-up'
$ws.Cells.Item(62, 2).Value = 'Theo dõi'
$ws.Cells.Item(62, 3).ClearContents()
$ws.Cells.Item(63, 1).Value = 'from cultures that have growth'
$ws.Cells.Item(63, 2).Value = 'từ mẫu cấy máu'
$ws.Cells.Item(63, 3).ClearContents()
$ws.Cells.Item(64, 1).Value = 'Generate .acorn file'
$ws.Cells.Item(64, 2).Value = 'Tạo tệp .acorn'
$ws.Cells.Item(64, 3).ClearContents()
$ws.Cells.Item(65, 1).Value = 'Generate and load .acorn from clinical and lab data'
$ws.Cells.Item(65, 2).Value = 'TBT'
$ws.Cells.Item(65, 3).Value = 'new'
$ws.Cells.Item(66, 1).Value = 'Generating .acorn'
$ws.Cells.Item(66, 2).Value = 'Đang tạo .acorn'
$ws.Cells.Item(66, 3).ClearContents()
$ws.Cells.Item(67, 1).Value = 'Get data from REDCap'
$ws.Cells.Item(67, 2).Value = 'Lấy dữ liệu từ REDCap'
$ws.Cells.Item(67, 3).ClearContents()
$ws.Cells.Item(68, 1).Value = 'Get the latest production release'
$ws.Cells.Item(68, 2).Value = 'TBT'
$ws.Cells.Item(68, 3).Value = 'new'
$ws.Cells.Item(69, 1).Value = 'Growth / No Growth'
$ws.Cells.Item(69, 2).Value = 'Cấy ra vi khuẩn/ Cấy không ra vi khuẩn'
$ws.Cells.Item(69, 3).ClearContents()
$ws.Cells.Item(70, 1).Value = 'HAI point prevalence by '
$ws.Cells.Item(70, 2).Value = 'TBT'
$ws.Cells.Item(70, 3).Value = 'new'
$ws.Cells.Item(71, 1).Value = 'HAI Prevalence'
$ws.Cells.Item(71, 2).Value = 'Tỷ lệ hiện mắc HAI'
$ws.Cells.Item(71, 3).ClearContents()
$ws.Cells.Item(72, 1).Value = 'Horizontal bars show the size of a set of SR results while vertical bars show the number of resistant isolates for the corresponding antibiotic.'
$ws.Cells.Item(72, 2).Value = 'TBT'
$ws.Cells.Item(72, 3).Value = 'new'
$ws.Cells.Item(73, 1).Value = 'Info on loaded .acorn'
$ws.Cells.Item(73, 2).Value = 'TBT'
$ws.Cells.Item(73, 3).Value = 'new'
$ws.Cells.Item(74, 1).Value = 'Initial &amp; Final Surveillance Diagnosis'
$ws.Cells.Item(74, 2).Value = 'Chẩn đoán giám sát ban đầu và cuối cùng'
$ws.Cells.Item(74, 3).ClearContents()
$ws.Cells.Item(75, 1).Value = 'Isolates'
$ws.Cells.Item(75, 2).Value = 'Số vi khuẩn cấy được'
$ws.Cells.Item(75, 3).ClearContents()
$ws.Cells.Item(76, 1).Value = 'Issue detected with REDCap data. Please report to ACORN data managers. Until resolution, only existing .acorn files can be used.'
$ws.Cells.Item(76, 2).Value = 'TBT'
$ws.Cells.Item(76, 3).Value = 'new'
$ws.Cells.Item(77, 1).Value = 'It might take a couple of minutes. This window will close on completion.'
$ws.Cells.Item(77, 2).Value = 'Thao tác có thể mất vài phút. Cửa sổ này sẽ đóng lại sau khi hoàn thành.'
$ws.Cells.Item(77, 3).ClearContents()
$ws.Cells.Item(78, 1).Value = 'Lab data successfully processed!'
$ws.Cells.Item(78, 2).Value = 'Dữ liệu phòng thí nghiệm đã được xử lý thành công!'
$ws.Cells.Item(78, 3).ClearContents()
$ws.Cells.Item(79, 1).Value = 'Lab data successfully provided.'
$ws.Cells.Item(79, 2).Value = 'Dữ liệu phòng xét nghiệm được cung cấp thành công'
$ws.Cells.Item(79, 3).ClearContents()
$ws.Cells.Item(80, 1).Value = 'Lab dataset contains the minimal columns.'
$ws.Cells.Item(80, 2).Value = 'Tập dữ liệu phòng thí nghiệm chứa các cột tối thiểu.'
$ws.Cells.Item(80, 3).ClearContents()
$ws.Cells.Item(81, 1).Value = 'Lab dataset does not contains the minimal columns.'
$ws.Cells.Item(81, 2).Value = 'Tập dữ liệu phòng thí nghiệm không chứa các cột tối thiểu.'
$ws.Cells.Item(81, 3).ClearContents()
$ws.Cells.Item(82, 1).Value = 'Language'
$ws.Cells.Item(82, 2).Value = 'Ngôn ngữ'
$ws.Cells.Item(82, 3).ClearContents()
$ws.Cells.Item(83, 1).Value = 'Load .acorn'
$ws.Cells.Item(83, 2).Value = 'Đang tải tệp .acorn'
$ws.Cells.Item(83, 3).ClearContents()
$ws.Cells.Item(84, 1).Value = 'Load .acorn from cloud'
$ws.Cells.Item(84, 2).Value = 'TBT'
$ws.Cells.Item(84, 3).Value = 'new'
$ws.Cells.Item(85, 1).Value = 'Load .acorn from local file'
$ws.Cells.Item(85, 2).Value = 'TBT'
$ws.Cells.Item(85, 3).Value = 'new'
$ws.Cells.Item(86, 1).Value = 'Load selected .acorn'
$ws.Cells.Item(86, 2).Value = 'Đang tải tệp .acorn được chọn'
$ws.Cells.Item(86, 3).ClearContents()
$ws.Cells.Item(87, 1).Value = 'Loading data.'
$ws.Cells.Item(87, 2).Value = 'Đang tải dữ liệu.'
$ws.Cells.Item(87, 3).ClearContents()
$ws.Cells.Item(88, 1).Value = 'Log in'
$ws.Cells.Item(88, 2).Value = 'Đăng nhập'
$ws.Cells.Item(88, 3).ClearContents()
$ws.Cells.Item(89, 1).Value = 'Microbiology'
$ws.Cells.Item(89, 2).Value = 'Vi sinh vật học'
$ws.Cells.Item(89, 3).ClearContents()
$ws.Cells.Item(90, 1).Value = 'Most frequent 10 organisms in the plot and complete listing in the table. Contaminants are in red.'
$ws.Cells.Item(90, 2).Value = '10 loại vi sinh vậy thường gặp nhất và liệt kê đầy đủ trong bảng. Tạp nhiễm có màu đỏ.'
$ws.Cells.Item(90, 3).ClearContents()
$ws.Cells.Item(91, 1).Value = 'No .acorn data loaded.'
$ws.Cells.Item(91, 2).Value = 'Không có dữ liệu .acorn nào được tải.'
$ws.Cells.Item(91, 3).ClearContents()
$ws.Cells.Item(92, 1).Value = 'No Blood Culture'
$ws.Cells.Item(92, 2).Value = 'Không cấy máu'
$ws.Cells.Item(92, 3).ClearContents()
$ws.Cells.Item(93, 1).Value = 'Not connected to internet.'
$ws.Cells.Item(93, 2).Value = 'Không có kết nối internet.'
$ws.Cells.Item(93, 3).ClearContents()
$ws.Cells.Item(94, 1).Value = 'Number of specimens per specimen type'
$ws.Cells.Item(94, 2).Value = 'Số lượng mẫu cho mỗi loại mẫu'
$ws.Cells.Item(94, 3).ClearContents()
$ws.Cells.Item(95, 1).Value = 'Occupancy rate per type of ward per month'
$ws.Cells.Item(95, 2).Value = 'Tỷ lệ chiếm chỗ của loại khoa mỗi tháng'
$ws.Cells.Item(95, 3).ClearContents()
$ws.Cells.Item(96, 1).Value = 'of blood cultures grew a potential contaminant.'
$ws.Cells.Item(96, 2).Value = 'Cấy máu phát triển một tạp nhiễm tiềm năng'
$ws.Cells.Item(96, 3).ClearContents()
$ws.Cells.Item(97, 1).Value = 'of cultures have growth.'
$ws.Cells.Item(97, 2).Value = 'Các mẫu cấy ra vi khuẩn'
$ws.Cells.Item(97, 3).ClearContents()
$ws.Cells.Item(98, 1).Value = 'of enrolments with blood culture.'
$ws.Cells.Item(98, 2).Value = 'thu tuyển với cấy máu'
$ws.Cells.Item(98, 3).ClearContents()
$ws.Cells.Item(99, 1).Value = 'of Target Pathogens'
$ws.Cells.Item(99, 2).Value = 'tác nhân gây bệnh mục tiêu'
$ws.Cells.Item(99, 3).ClearContents()
$ws.Cells.Item(100, 1).Value = 'Only isolates that have been tested against all of the drugs are included in the upset plot.'
$ws.Cells.Item(100, 2).Value = 'TBT'
$ws.Cells.Item(100, 3).Value = 'new'
$ws.Cells.Item(101, 1).Value = 'Overview'
$ws.Cells.Item(101, 2).Value = 'Tổng quan'
$ws.Cells.Item(101, 3).ClearContents()
$ws.Cells.Item(102, 1).Value = 'Password'
$ws.Cells.Item(102, 2).Value = 'Mật khẩu'
$ws.Cells.Item(102, 3).ClearContents()
$ws.Cells.Item(103, 1).Value = 'Patient Age Distribution'
$ws.Cells.Item(103, 2).Value = 'Phân bố tuổi bệnh nhân'
$ws.Cells.Item(103, 3).ClearContents()
$ws.Cells.Item(104, 1).Value = 'Patient Comorbidities'
$ws.Cells.Item(104, 2).Value = 'Bệnh nền của bệnh nhân'
$ws.Cells.Item(104, 3).ClearContents()
$ws.Cells.Item(105, 1).Value = 'Patient enrolments'
$ws.Cells.Item(105, 2).Value = 'Thu tuyển bệnh nhân'
$ws.Cells.Item(105, 3).ClearContents()
$ws.Cells.Item(106, 1).Value = 'Patients Transferred'
$ws.Cells.Item(106, 2).Value = 'Bệnh nhân được chuyển'
$ws.Cells.Item(106, 3).ClearContents()
$ws.Cells.Item(107, 1).Value = 'Please log in'
$ws.Cells.Item(107, 2).Value = 'Hãy đăng nhập'
$ws.Cells.Item(107, 3).ClearContents()
$ws.Cells.Item(108, 1).Value = 'Processing lab data.'
$ws.Cells.Item(108, 2).Value = 'Xử lý dữ liệu phòng thí nghiệm.'
$ws.Cells.Item(108, 3).ClearContents()
$ws.Cells.Item(109, 1).Value = 'Reading lab data.'
$ws.Cells.Item(109, 2).Value = 'Đọc dữ liệu phòng thí nghiệm.'
$ws.Cells.Item(109, 3).ClearContents()
$ws.Cells.Item(110, 1).Value = 'Remove ''Not Cultured'' specimens'
$ws.Cells.Item(110, 2).Value = 'TBT'
$ws.Cells.Item(110, 3).Value = 'new'
$ws.Cells.Item(111, 1).Value = 'Remove blood culture contaminants from the following visualizations'
$ws.Cells.Item(111, 2).Value = 'Loại bỏ các tạp nhiễm trong quá trình cấy máu khỏi các hình ảnh trực quan'
$ws.Cells.Item(111, 3).ClearContents()
$ws.Cells.Item(112, 1).Value = 'Reset Enrolments Filters'
$ws.Cells.Item(112, 2).Value = 'Cài đặt lại bộ lọc thu tuyển'
$ws.Cells.Item(112, 3).ClearContents()
$ws.Cells.Item(113, 1).Value = 'Resistance to 3rd gen. Cephalosporins Over Time'
$ws.Cells.Item(113, 2).Value = 'Kháng Cephalosporins thế hệ 3 theo thời gian.'
$ws.Cells.Item(113, 3).ClearContents()
$ws.Cells.Item(114, 1).Value = 'Resistance to Carbapenems Over Time'
$ws.Cells.Item(114, 2).Value = 'Kháng Carbapenems theo thời gian'
$ws.Cells.Item(114, 3).ClearContents()
$ws.Cells.Item(115, 1).Value = 'Resistance to Fluoroquinolones Over Time'
$ws.Cells.Item(115, 2).Value = 'Kháng Fluoroquinolones theo thời gian'
$ws.Cells.Item(115, 3).ClearContents()
$ws.Cells.Item(116, 1).Value = 'Resistance to Oxacillin Over Time'
$ws.Cells.Item(116, 2).Value = 'Kháng Oxacillin theo thời gian'
$ws.Cells.Item(116, 3).ClearContents()
$ws.Cells.Item(117, 1).Value = 'Resistance to Penicillin G - meningitis Over Time'
$ws.Cells.Item(117, 2).Value = 'Kháng Penicillin G -meningitis theo thời gian'
$ws.Cells.Item(117, 3).ClearContents()
$ws.Cells.Item(118, 1).Value = 'Resistance to Penicillin G Over Time'
$ws.Cells.Item(118, 2).Value = 'Kháng Penicillin G theo thời gian'
$ws.Cells.Item(118, 3).ClearContents()
$ws.Cells.Item(119, 1).Value = 'Retriving data from REDCap server.'
$ws.Cells.Item(119, 2).Value = 'Truy xuất dữ liệu từ máy chủ REDCap.'
$ws.Cells.Item(119, 3).ClearContents()
$ws.Cells.Item(120, 1).Value = 'Save .acorn file'
$ws.Cells.Item(120, 2).Value = 'Lưu tệp .acorn'
$ws.Cells.Item(120, 3).ClearContents()
$ws.Cells.Item(121, 1).Value = 'Save acorn data'
$ws.Cells.Item(121, 2).Value = 'Lưu dữ liệu acorn'
$ws.Cells.Item(121, 3).ClearContents()
$ws.Cells.Item(122, 1).Value = 'Save on Server'
$ws.Cells.Item(122, 2).Value = 'Lưu trên máy chủ'
$ws.Cells.Item(122, 3).ClearContents()
$ws.Cells.Item(123, 1).Value = 'See Breakdown by Ward'
$ws.Cells.Item(123, 2).Value = 'Xem Phân tích theo Khoa'
$ws.Cells.Item(123, 3).ClearContents()
$ws.Cells.Item(124, 1).Value = 'See by Week'
$ws.Cells.Item(124, 2).Value = 'Xem theo Tuần'
$ws.Cells.Item(124, 3).ClearContents()
$ws.Cells.Item(125, 1).Value = 'Show antibiotics combinations'
$ws.Cells.Item(125, 2).Value = 'TBT'
$ws.Cells.Item(125, 3).Value = 'new'
$ws.Cells.Item(126, 1).Value = 'Show comorbidities combinations'
$ws.Cells.Item(126, 2).Value = 'Hiển thị các bệnh nền'
$ws.Cells.Item(126, 3).ClearContents()
$ws.Cells.Item(127, 1).Value = 'SIR Evaluation'
$ws.Cells.Item(127, 2).Value = 'TBT'
$ws.Cells.Item(127, 3).Value = 'new'
$ws.Cells.Item(128, 1).Value = 'Some D28 records (F04) don''t have a matching patient enrolment (F01).'
$ws.Cells.Item(128, 2).Value = 'Một số Phiếu dữ liệu ngày 28 (F04) không có một bệnh nhân thu tuyển (F01) tương ứng'
$ws.Cells.Item(128, 3).ClearContents()
$ws.Cells.Item(129, 1).Value = 'Some dates of enrolment for HAI patients do have a matching date in the HAI survey dataset'
$ws.Cells.Item(129, 2).Value = 'Một số ngày thu tuyển của bẹnh nhân HAI  không tương ứng với ngày thực hiện giám sát điểm đã được định sẵn'
$ws.Cells.Item(129, 3).ClearContents()
$ws.Cells.Item(130, 1).Value = 'Some hospital outcome records (F03) don''t have a matching infection episode (F02). These records have been removed.'
$ws.Cells.Item(130, 2).Value = 'Một số Phiếu ra viện (F03) không có phiếu đợt nhiễm trùng tương ứng (F02). Những phiếu này đã bị xóa'
$ws.Cells.Item(130, 3).ClearContents()
$ws.Cells.Item(131, 1).Value = 'Some hospital outcome records (F03) don''t have a matching patient enrolment (F01).'
$ws.Cells.Item(131, 2).Value = 'Mỗi phiếu xuất viện (F03) không có Phiếu thu tuyển (F01) tương ứng'
$ws.Cells.Item(131, 3).ClearContents()
$ws.Cells.Item(132, 1).Value = 'Some infection episode records (F02) don''t have a matching patient enrolment (F01). These records have been removed.'
$ws.Cells.Item(132, 2).Value = 'Một số Phiếu đợt nhiễm trùng (F02) không có phiếu thu tuyển (F01) tương ứng. Những phiếu này đã bị xóa'
$ws.Cells.Item(132, 3).ClearContents()
$ws.Cells.Item(133, 1).Value = 'Some records with a missing ACORN ID. These records have been removed.'
$ws.Cells.Item(133, 2).Value = 'TBT'
$ws.Cells.Item(133, 3).Value = 'new'
$ws.Cells.Item(134, 1).Value = 'Specimen Types'
$ws.Cells.Item(134, 2).Value = 'Loại mẫu'
$ws.Cells.Item(134, 3).ClearContents()
$ws.Cells.Item(135, 1).Value = 'Specimens'
$ws.Cells.Item(135, 2).Value = 'TBT'
$ws.Cells.Item(135, 3).Value = 'new'
$ws.Cells.Item(136, 1).Value = 'Specimens Collected'
$ws.Cells.Item(136, 2).Value = 'Mẫu bệnh phẩm thu thập được'
$ws.Cells.Item(136, 3).ClearContents()
$ws.Cells.Item(137, 1).Value = 'specimens per enrolment'
$ws.Cells.Item(137, 2).Value = 'Số mẫu bệnh phẩm/ số thu tuyển'
$ws.Cells.Item(137, 3).ClearContents()
$ws.Cells.Item(138, 1).Value = 'Successfully combined clinical and lab data into .acorn file'
$ws.Cells.Item(138, 2).Value = 'Kết hợp thành công dữ liệu lâm sàng và dữ liệu phòng xét nghiệm vào file .acorn'
$ws.Cells.Item(138, 3).ClearContents()
$ws.Cells.Item(139, 1).Value = 'Successfully loaded data.'
$ws.Cells.Item(139, 2).Value = 'Tải dữ liệu thành công.'
$ws.Cells.Item(139, 3).ClearContents()
$ws.Cells.Item(140, 1).Value = 'Successfully logged in.'
$ws.Cells.Item(140, 2).Value = 'Đăng nhập thành công.'
$ws.Cells.Item(140, 3).ClearContents()
$ws.Cells.Item(141, 1).Value = 'Successfully saved .acorn file in the cloud. You can now explore acorn data.'
$ws.Cells.Item(141, 2).Value = 'Đã lưu thành công tệp .acorn trên đám mây. Bây giờ bạn có thể khám phá dữ liệu acorn.'
$ws.Cells.Item(141, 3).ClearContents()
$ws.Cells.Item(142, 1).Value = 'Successfully saved .acorn file locally.'
$ws.Cells.Item(142, 2).Value = 'Đã lưu thành công tệp .acorn tại điểm nghiên cứu.'
$ws.Cells.Item(142, 3).ClearContents()
$ws.Cells.Item(143, 1).Value = 'Supply first valid clinical and lab data.'
$ws.Cells.Item(143, 2).Value = 'Cung cấp dữ liệu lâm sàng và dữ liệu phòng thí nghiệm hợp lệ đầu tiên.'
$ws.Cells.Item(143, 3).ClearContents()
$ws.Cells.Item(144, 1).Value = 'Susceptible &amp; Intermediate are always combined in this visualisation of co-resistances.'
$ws.Cells.Item(144, 2).Value = 'TBT'
$ws.Cells.Item(144, 3).Value = 'new'
$ws.Cells.Item(145, 1).Value = 'The 10 most common initial-final diagnosis combinations:'
$ws.Cells.Item(145, 2).Value = '10 kết hợp chẩn đoán ban đầu-cuối cùng phổ biến nhất:'
$ws.Cells.Item(145, 3).ClearContents()
$ws.Cells.Item(146, 1).Value = 'The following ''patient id'' are atypical cases (one HCAI/CAI with early HAI but no overlap):'
$ws.Cells.Item(146, 2).Value = 'TBT'
$ws.Cells.Item(146, 3).Value = 'new'
$ws.Cells.Item(147, 1).Value = 'The following ''patient id'' are problem case (overlapping specimen collection windows):'
$ws.Cells.Item(147, 2).Value = 'ID bệnh nhân'' sau là vấn đề(các cửa sổ thu thập mẫu bệnh phẩm chồng chéo):'
$ws.Cells.Item(147, 3).ClearContents()
$ws.Cells.Item(148, 1).Value = 'The REDCap dataset is empty/in wrong format. Please contact ACORN support.'
$ws.Cells.Item(148, 2).Value = 'TBT'
$ws.Cells.Item(148, 3).Value = 'new'
$ws.Cells.Item(149, 1).Value = 'The REDCap dataset is in the right format.'
$ws.Cells.Item(149, 2).Value = 'TBT'
$ws.Cells.Item(149, 3).Value = 'new'
$ws.Cells.Item(150, 1).Value = 'There are D28 follow-up done before the expected D28 date.'
$ws.Cells.Item(150, 2).Value = 'TBT'
$ws.Cells.Item(150, 3).Value = 'new'
$ws.Cells.Item(151, 1).Value = 'There are multiple F02 with identical ACORN ID, admission date, and episode enrolment date.'
$ws.Cells.Item(151, 2).Value = 'TBT'
$ws.Cells.Item(151, 3).Value = 'new'
$ws.Cells.Item(152, 1).Value = 'There are no atypical case (one HCAI/CAI with early HAI but no overlap).'
$ws.Cells.Item(152, 2).Value = 'TBT'
$ws.Cells.Item(152, 3).Value = 'new'
$ws.Cells.Item(153, 1).Value = 'There are no D28 follow-up done before the expected D28 date.'
$ws.Cells.Item(153, 2).Value = 'TBT'
$ws.Cells.Item(153, 3).Value = 'new'
$ws.Cells.Item(154, 1).Value = 'There are no isolate with valid AST results. Please contact ACORN support.'
$ws.Cells.Item(154, 2).Value = 'Không có chủng vi khuẩn với kết quả Kháng sinh đồ có giá trị. Vui lòng liên hệ bộ phận hỗ trợ ACORN'
$ws.Cells.Item(154, 3).ClearContents()
$ws.Cells.Item(155, 1).Value = 'There are no multiple F02 with identical ACORN ID, admission date, and episode enrolment date.'
$ws.Cells.Item(155, 2).Value = 'TBT'
$ws.Cells.Item(155, 3).Value = 'new'
$ws.Cells.Item(156, 1).Value = 'There are no problem case (overlapping specimen collection windows)'
$ws.Cells.Item(156, 2).Value = 'Không có vấn đề (các cửa sổ thu thập mẫu bệnh phẩm chồng chéo)'
$ws.Cells.Item(156, 3).ClearContents()
$ws.Cells.Item(157, 1).Value = 'There are rows for which ''specdate'' are after today.'
$ws.Cells.Item(157, 2).Value = 'Có dòng mà "specdate" sau ngày hôm nay'
$ws.Cells.Item(157, 3).ClearContents()
$ws.Cells.Item(158, 1).Value = 'There are rows with missing ''orgname''.'
$ws.Cells.Item(158, 2).Value = 'Có dòng thiếu "orgname"'
$ws.Cells.Item(158, 3).ClearContents()
$ws.Cells.Item(159, 1).Value = 'There are rows with missing ''patid''.'
$ws.Cells.Item(159, 2).Value = 'Có dòng thiếu "patid"'
$ws.Cells.Item(159, 3).ClearContents()
$ws.Cells.Item(160, 1).Value = 'There are rows with missing ''specdate''.'
$ws.Cells.Item(160, 2).Value = 'Có dòng thiếu "specdate"'
$ws.Cells.Item(160, 3).ClearContents()
$ws.Cells.Item(161, 1).Value = 'There are rows with missing ''specgroup''.'
$ws.Cells.Item(161, 2).Value = 'Có dòng thiếu "specgroup"'
$ws.Cells.Item(161, 3).ClearContents()
$ws.Cells.Item(162, 1).Value = 'There are rows with missing ''specid''.'
$ws.Cells.Item(162, 2).Value = 'Có dòng thiếu "specid"'
$ws.Cells.Item(162, 3).ClearContents()
$ws.Cells.Item(163, 1).Value = 'There is a critical issue with clinical data. The issue should be fixed in REDCap.'
$ws.Cells.Item(163, 2).Value = 'Có một vấn đề nghiêm trọng với dữ liệu lâm sàng. Vấn đề này nên được chỉnh sửa trên REDCap.'
$ws.Cells.Item(163, 3).ClearContents()
$ws.Cells.Item(164, 1).Value = 'There is no data to display for this organism.'
$ws.Cells.Item(164, 2).Value = 'Không có dữ liệu mô tả cho vi sinh vật này'
$ws.Cells.Item(164, 3).ClearContents()
$ws.Cells.Item(165, 1).Value = 'There is no HAI survey data'
$ws.Cells.Item(165, 2).Value = 'Không có dữ liệu khảo sát HAI'
$ws.Cells.Item(165, 3).ClearContents()
$ws.Cells.Item(166, 1).Value = 'Trying to save .acorn file on server.'
$ws.Cells.Item(166, 2).Value = 'Đang cố gắng lưu tệp .acorn trên máy chủ.'
$ws.Cells.Item(166, 3).ClearContents()
$ws.Cells.Item(167, 1).Value = 'Updated Charlson Comorbidity Index (uCCI)'
$ws.Cells.Item(167, 2).Value = 'TBT'
$ws.Cells.Item(167, 3).Value = 'new'
$ws.Cells.Item(168, 1).Value = 'User'
$ws.Cells.Item(168, 2).Value = 'Người sử dụng'
$ws.Cells.Item(168, 3).ClearContents()
$ws.Cells.Item(169, 1).Value = 'Variables in Table:'
$ws.Cells.Item(169, 2).Value = 'Các biến số trong Bảng:'
$ws.Cells.Item(169, 3).ClearContents()
$ws.Cells.Item(170, 1).Value = 'Ward Occupancy Rates'
$ws.Cells.Item(170, 2).Value = 'Tỷ lệ số giường bệnh của khoa được lấp đầy'
$ws.Cells.Item(170, 3).ClearContents()
$ws.Cells.Item(171, 1).Value = 'We couldn''t download the lab codes file. Please contact ACORN support.'
$ws.Cells.Item(171, 2).Value = 'Chúng tôi không thể tải xuống tệp code phòng xét nghiệm. Vui lòng liên hệ bộ phận hỗ trợ của ACORN'
$ws.Cells.Item(171, 3).ClearContents()
$ws.Cells.Item(172, 1).Value = 'We couldn''t download the lab data dictionary. Please contact ACORN support'
$ws.Cells.Item(172, 2).Value = 'TBT'
$ws.Cells.Item(172, 3).Value = 'new'
$ws.Cells.Item(173, 1).Value = 'Welcome'
$ws.Cells.Item(173, 2).Value = 'Xin chào'
$ws.Cells.Item(173, 3).ClearContents()
$ws.Cells.Item(174, 1).Value = 'What do you want to do?'
$ws.Cells.Item(174, 2).Value = 'Bạn muốn làm gì?'
$ws.Cells.Item(174, 3).ClearContents()
$ws.Cells.Item(175, 1).Value = 'With Microbiology'
$ws.Cells.Item(175, 2).Value = 'Với vi sinh vật học'
$ws.Cells.Item(175, 3).ClearContents()
$ws.Cells.Item(176, 1).Value = 'Wrong connection credentials.'
$ws.Cells.Item(176, 2).Value = 'Thông tin đăng nhập kết nối sai'
$ws.Cells.Item(176, 3).ClearContents()
$ws.Cells.Item(177, 1).Value = 'You are running ACORN dashboard'
$ws.Cells.Item(177, 2).Value = 'TBT'
$ws.Cells.Item(177, 3).Value = 'new'
$ws.Cells.Item(178, 1).Value = 'You can check here if it''s the latest production release.'
$ws.Cells.Item(178, 2).Value = 'TBT'
$ws.Cells.Item(178, 3).Value = 'new'
$ws.Cells.Item(179, 1).Value = 'Your ACORN dashboard is up to date'
$ws.Cells.Item(179, 2).Value = 'TBT'
$ws.Cells.Item(179, 3).Value = 'new'
$ws.Cells.Item(180, 1).Value = 'Follow us on Twitter'
$ws.Cells.Item(180, 2).Value = 'Theo dõi trên Twitter'
$ws.Cells.Item(180, 3).ClearContents()
$ws.Cells.Item(181, 1).Value = 'Records in Lab data and BSI forms:'
$ws.Cells.Item(181, 2).Value = 'Ghi chép ở dữ liệu phòng xét nghiệm và Phiếu BSI'
$ws.Cells.Item(181, 3).ClearContents()
$ws.Cells.Item(182, 1).Value = 'What is ACORN?'
$ws.Cells.Item(182, 2).Value = 'Nghiên cứu ACORN là gì?'
$ws.Cells.Item(182, 3).Value = 'deleted'
$ws.Cells.Item(183, 1).Value = 'About'
$ws.Cells.Item(183, 2).Value = 'Về chúng tôi'
$ws.Cells.Item(183, 3).Value = 'deleted'
$ws.Cells.Item(184, 1).Value = 'Site'
$ws.Cells.Item(184, 2).Value = 'Địa điểm'
$ws.Cells.Item(184, 3).Value = 'deleted'
$ws.Cells.Item(185, 1).Value = 'To log out, close the app.'
$ws.Cells.Item(185, 2).Value = 'Để thoát ra, vui lòng đóng ứng dụng.'
$ws.Cells.Item(185, 3).Value = 'deleted'
$ws.Cells.Item(186, 1).Value = 'upload a local acorn file.'
$ws.Cells.Item(186, 2).Value = 'Tải lên một tệp acorn của điểm nghiên cứu.'
$ws.Cells.Item(186, 3).Value = 'deleted'
$ws.Cells.Item(187, 1).Value = 'The REDCap dataset contains data.'
$ws.Cells.Item(187, 2).Value = 'Bộ dữ liệu REDCap có chứa dữ liệu.'
$ws.Cells.Item(187, 3).Value = 'deleted'
$ws.Cells.Item(188, 1).Value = 'The REDCap dataset is empty. Please contact ACORN support.'
$ws.Cells.Item(188, 2).Value = 'Bộ dữ liệu REDCap trống. Vui lòng liên hệ bộ phận hỗ trợ của ACORN.'
$ws.Cells.Item(188, 3).Value = 'deleted'
$ws.Cells.Item(189, 1).Value = 'The REDCap dataset column names do not match. Please contact ACORN support.'
$ws.Cells.Item(189, 2).Value = 'Tên các cột của bộ dữ liệu REDCap không khớp. Vui lòng liên hệ bộ phận hỗ trợ của ACORN.'
$ws.Cells.Item(189, 3).Value = 'deleted'
$ws.Cells.Item(190, 1).Value = 'The REDCap dataset column names match.'
$ws.Cells.Item(190, 2).Value = 'Tên các cột của bộ dữ liệu REDCap trùng khớp.'
$ws.Cells.Item(190, 3).Value = 'deleted'
$ws.Cells.Item(191, 1).Value = 'Clinical data successfully provided.'
$ws.Cells.Item(191, 2).Value = 'Dữ liệu lâm sàng đã được cung cấp thành công.'
$ws.Cells.Item(191, 3).Value = 'deleted'
$ws.Cells.Item(192, 1).Value = 'Clinical data not provided'
$ws.Cells.Item(192, 2).Value = 'Dữ liệu lâm sàng không được cung cấp'
$ws.Cells.Item(192, 3).Value = 'deleted'
$ws.Cells.Item(193, 1).Value = 'Lab data not provided'
$ws.Cells.Item(193, 2).Value = 'Dữ liệu phòng thí nghiệm không được cung cấp'
$ws.Cells.Item(193, 3).Value = 'deleted'
$ws.Cells.Item(194, 1).Value = 'No .acorn has been generated'
$ws.Cells.Item(194, 2).Value = 'Không có tệp .acorn được tạo'
$ws.Cells.Item(194, 3).Value = 'deleted'
$ws.Cells.Item(195, 1).Value = 'No .acorn has been saved'
$ws.Cells.Item(195, 2).Value = 'Không có tệp .acorn được lưu'
$ws.Cells.Item(195, 3).Value = 'deleted'
$ws.Cells.Item(196, 1).Value = 'Generate .acorn from clinical and lab data'
$ws.Cells.Item(196, 2).Value = 'Tạo tệp .acorn từ dữ liệu lâm sàng và dữ liệu phòng thí nghiệm'
$ws.Cells.Item(196, 3).Value = 'deleted'
$ws.Cells.Item(197, 1).Value = 'Load existing .acorn from cloud'
$ws.Cells.Item(197, 2).Value = 'Tải tệp .acorn sẵn có từ đám mây'
$ws.Cells.Item(197, 3).Value = 'deleted'
$ws.Cells.Item(198, 1).Value = 'Load existing .acorn from local file'
$ws.Cells.Item(198, 2).Value = 'Tải tệp .acorn sẵn có từ tệp tại điểm nghiên cứu'
$ws.Cells.Item(198, 3).Value = 'deleted'
$ws.Cells.Item(199, 1).Value = 'Critical issue detected: no data or wrong data format on REDCap server. Please report to ACORN data managers. Until resolution, only existing .acorn files can be used.'
$ws.Cells.Item(199, 2).Value = 'Phát hiện vấn đề nghiêm trọng: không có dữ liệu hoặc định dạng dữ liệu sai trên máy chủ REDCap. Vui lòng báo cáo cho quản lý dữ liệu của ACORN. Cho tới khi tìm ra giải pháp, vui lòng chỉ sử dụng các tệp .acorn hiện có.'
$ws.Cells.Item(199, 3).Value = 'deleted'
$ws.Cells.Item(200, 1).Value = 'Successfully saved .acorn file locally. You can now explore acorn data.'
$ws.Cells.Item(200, 2).Value = 'Đã lưu thành công tệp .acorn tại điểm nghiên cứu. Bây giờ bạn có thể tìm hiểu dữ liệu acorn.'
$ws.Cells.Item(200, 3).Value = 'deleted'
$ws.Cells.Item(201, 1).Value = 'Select lab data format:'
$ws.Cells.Item(201, 2).Value = 'Lựa chọn định dạng dữ liệu phòng thí nghiệm:'
$ws.Cells.Item(201, 3).Value = 'deleted'
$ws.Cells.Item(202, 1).Value = 'Specimens, Isolates'
$ws.Cells.Item(202, 2).Value = 'Số mẫu, Chủng phân lập được'
$ws.Cells.Item(202, 3).Value = 'deleted'
$ws.Cells.Item(203, 1).Value = 'Problem with credentials. Please contact ACORN support.'
$ws.Cells.Item(203, 2).Value = 'Vấn đề với thông tin xác thực. Vui lòng liên hệ với bộ phận hỗ trợ của ACORN.'
$ws.Cells.Item(203, 3).Value = 'deleted'
$ws.Cells.Item(204, 1).Value = 'Blood Culture'
$ws.Cells.Item(204, 2).Value = 'Mẫu cấy máu'
$ws.Cells.Item(204, 3).Value = 'deleted'
$ws.Cells.Item(205, 1).Value = 'Other Specimens:'
$ws.Cells.Item(205, 2).Value = 'Các mẫu khác'
$ws.Cells.Item(205, 3).Value = 'deleted'
$ws.Cells.Item(206, 1).Value = 'No deduplication of isolates'
$ws.Cells.Item(206, 2).Value = 'Không trùng lặp các mẫu phân lập'
$ws.Cells.Item(206, 3).Value = 'deleted'
$ws.Cells.Item(207, 1).Value = 'Deduplication by patient-episode'
$ws.Cells.Item(207, 2).Value = 'Sự trùng lặp đợt nhiễm trùng của bệnh nhân'
$ws.Cells.Item(207, 3).Value = 'deleted'
$ws.Cells.Item(208, 1).Value = 'Deduplication by patient ID'
$ws.Cells.Item(208, 2).Value = 'Sự trùng lặp ID của bệnh nhân'
$ws.Cells.Item(208, 3).Value = 'deleted'
$ws.Cells.Item(209, 1).Value = 'Use heuristic for time unit'
$ws.Cells.Item(209, 2).Value = 'Sử dụng heuristic cho đơn vị thời gian'
$ws.Cells.Item(209, 3).Value = 'deleted'
$ws.Cells.Item(210, 1).Value = 'Display by month'
$ws.Cells.Item(210, 2).Value = 'Hiển thị theo tháng'
$ws.Cells.Item(210, 3).Value = 'deleted'
$ws.Cells.Item(211, 1).Value = 'Display by year'
$ws.Cells.Item(211, 2).Value = 'Hiển thị theo năm'
$ws.Cells.Item(211, 3).Value = 'deleted'
$ws.Cells.Item(212, 1).Value = 'HAI point prevalence by type of ward'
$ws.Cells.Item(212, 2).Value = 'Tỷ lệ hiện mắc điểm HAI theo loại khoa'
$ws.Cells.Item(212, 3).Value = 'deleted'
$ws.Cells.Item(213, 1).Value = 'We couldn''t download the lab data dictionary. Please contact ACORN support.'
$ws.Cells.Item(213, 2).Value = 'Chúng tôi không thể tải xuống từ điển dữ liệu phòng xét nghiệm. Vui lòng liện hệ bộ phận hỗ trợ của ACORN'
$ws.Cells.Item(213, 3).Value = 'deleted'
$ws.Cells.Item(214, 1).Value = 'REDCap data could not be downloaded. Please try again.'
$ws.Cells.Item(214, 2).Value = 'Dữ liệu REDCap không thể tải xuống. Vui lòng thử lại'
$ws.Cells.Item(214, 3).Value = 'deleted'
$ws.Cells.Item(215, 1).Value = 'All records have an ACORN ID.'
$ws.Cells.Item(215, 2).Value = 'Tất cả các phiếu thu thập đều có một ACORN ID'
$ws.Cells.Item(215, 3).Value = 'deleted'
$ws.Cells.Item(216, 1).Value = 'Not all records have an ACORN ID.'
$ws.Cells.Item(216, 2).Value = 'Không phải tất cả các phiếu thu thập đều có một ACORN ID'
$ws.Cells.Item(216, 3).Value = 'deleted'
$ws.Cells.Item(217, 1).Value = 'There are no atypical case (one CAI / early HAI but no overlap).'
$ws.Cells.Item(217, 2).Value = 'Không có trường hợp điển hình ( một CAI/ HAI sớm nhưng không trùng lặp)'
$ws.Cells.Item(217, 3).Value = 'deleted'
$ws.Cells.Item(218, 1).Value = 'The following ''patient id'' are atypical cases (one CAI / early HAI but no overlap):'
$ws.Cells.Item(218, 2).Value = 'ID bệnh nhân'' sau là các trường hợp không điển hình (một CAI / HAI sớm nhưng không trùng lặp)'
$ws.Cells.Item(218, 3).Value = 'deleted'